# Refresh the cryptos table to the latest scraped snapshot.
# Most rows only change Price (D) / Volume(1h) (E); rows 48, 49 and 51
# additionally swap in a different coin (Coin name + Link too).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.467.61'
$ws.Range('E2').Value = '  -2.69%  '
$ws.Range('D3').Value = '1.984.65'
$ws.Range('E3').Value = '  -3.33%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'244.37"
$ws.Range('E5').Value = '  -3.35%  '
$ws.Range('E6').Value = '  -3.36%  '
$ws.Range('D7').Value = "'59.16"
$ws.Range('E7').Value = '  -10.70%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = "'0.377"
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('D10').Value = "'57.93"
$ws.Range('E10').Value = '  -2.92%  '
$ws.Range('D11').Value = "'0.0815"
$ws.Range('E11').Value = '  +5.76%  '
$ws.Range('E12').Value = '  -1.11%  '
$ws.Range('D13').Value = "'23.81"
$ws.Range('E13').Value = '  +5.07%  '
$ws.Range('E14').Value = '  -5.27%  '
$ws.Range('D15').Value = "'14.03"
$ws.Range('E15').Value = '  -5.92%  '
$ws.Range('D16').Value = '2.277.58'
$ws.Range('E16').Value = '  -3.20%  '
$ws.Range('D17').Value = "'5.48"
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('D18').Value = '1.977.61'
$ws.Range('E18').Value = '  -3.64%  '
$ws.Range('D19').Value = '36.430.81'
$ws.Range('E19').Value = '  -2.27%  '
$ws.Range('D20').Value = "'70.73"
$ws.Range('E20').Value = '  -4.20%  '
$ws.Range('E21').Value = '  -1.78%  '
$ws.Range('D22').Value = "'5.32"
$ws.Range('E22').Value = '  -2.30%  '
$ws.Range('D23').Value = "'234.62"
$ws.Range('E23').Value = '  -2.48%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').Value = "'2.60"
$ws.Range('E25').Value = '  -3.05%  '
$ws.Range('E26').Value = '  -4.32%  '
$ws.Range('D27').Value = "'10.22"
$ws.Range('E27').Value = '  +2.90%  '
$ws.Range('D28').Value = "'161.97"
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').Value = "'19.93"
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('D30').Value = "'0.131"
$ws.Range('E30').Value = '  +10.69%  '
$ws.Range('D31').Value = "'0.120"
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('D33').Value = "'4.92"
$ws.Range('E33').Value = '  -7.62%  '
$ws.Range('D34').Value = "'0.0631"
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('D35').Value = "'4.44"
$ws.Range('E35').Value = '  -6.14%  '
$ws.Range('D36').Value = "'6.31"
$ws.Range('E36').Value = '  +4.89%  '
$ws.Range('D37').Value = "'2.28"
$ws.Range('E37').Value = '  -7.63%  '
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('D39').Value = "'1.77"
$ws.Range('E39').Value = '  -4.51%  '
$ws.Range('D40').Value = "'3.07"
$ws.Range('E40').Value = '  +2.48%  '
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('E42').Value = '  -7.72%  '
$ws.Range('D43').Value = "'2.90"
$ws.Range('E43').Value = '  -3.79%  '
$ws.Range('E44').Value = '  -2.35%  '
$ws.Range('D45').Value = "'1.09"
$ws.Range('E45').Value = '  -4.97%  '
$ws.Range('D46').Value = "'92.92"
$ws.Range('E46').Value = '  -3.94%  '
$ws.Range('D47').Value = "'16.29"
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '1.385.42'
$ws.Range('E48').Value = '  -2.97%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = "'7.57"
$ws.Range('E49').Value = '  -5.34%  '
$ws.Range('E50').Value = '  -2.89%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = "'45.13"
$ws.Range('E51').Value = '  -3.90%  '
